$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 61, pushing all subsequent
# records (previously rows 61-141) down by one row (now rows 62-142).
$ws.Rows("61:61").Insert()

$ws.Range("A61").Value = 4
$ws.Range("B61").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C61").Value = "Los Lagos"
$ws.Range("D61").Value = 44579
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = 100112009
$ws.Range("G61").Value = "Acelga"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 80
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = 10000
$ws.Range("N61").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O61").Value = "Región de La Araucanía"
$ws.Range("P61").Value = 833
$ws.Range("Q61").Value = 12
$ws.Range("R61").Value = "Hortaliza"
